$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 25167
$ws.Range("J12").Value = 120
$ws.Range("L12").Value = 120
$ws.Range("N12").Value = -460
$ws.Range("H38").Value = 1422.8182
$ws.Range("I38").Value = 1115.6
$ws.Range("J38").Value = 4495
$ws.Range("K38").Value = 3346.8
$ws.Range("L38").Value = 13485
$ws.Range("M38").Value = -2974.8
$ws.Range("N38").Value = -14229
$ws.Range("H39").Value = 6806.4165
$ws.Range("I39").Value = 4289.5
$ws.Range("J39").Value = 7309.8
$ws.Range("K39").Value = 12868.5
$ws.Range("L39").Value = 21929.4
$ws.Range("M39").Value = -12572.5
$ws.Range("N39").Value = -22521.4
$ws.Range("H74").Value = 1126971.8
$ws.Range("I74").Value = 1126971.8
$ws.Range("K74").Value = 1126971.8
$ws.Range("M74").Value = -1126035.8
$ws.Range("H76").Value = 4158
$ws.Range("I76").Value = 3537
$ws.Range("K76").Value = 3537
$ws.Range("M76").Value = -3222
$ws.Range("H77").Value = 1126971.8
$ws.Range("I77").Value = 1126971.8
$ws.Range("K77").Value = 5634859
$ws.Range("M77").Value = -5630179
$ws.Range("H79").Value = 4158
$ws.Range("I79").Value = 3537
$ws.Range("K79").Value = 3537
$ws.Range("M79").Value = -2445
$ws.Range("H80").Value = 1311.6471
$ws.Range("I80").Value = 2066.5557
$ws.Range("J80").Value = 462.375
$ws.Range("K80").Value = 6199.6671
$ws.Range("L80").Value = 1387.125
$ws.Range("M80").Value = -5201.6671
$ws.Range("N80").Value = -3383.125
$ws.Range("H83").Value = 1311.6471
$ws.Range("I83").Value = 2066.5557
$ws.Range("J83").Value = 462.375
$ws.Range("K83").Value = 18599.0013
$ws.Range("L83").Value = 4161.375
$ws.Range("M83").Value = -13607.0013
$ws.Range("N83").Value = -14145.375
$ws.Range("H116").Value = 5000.4
$ws.Range("I116").Value = 5000.636
$ws.Range("J116").Value = 4998.6665
$ws.Range("K116").Value = 5000.636
$ws.Range("L116").Value = 4998.6665
$ws.Range("M116").Value = -1558.636
$ws.Range("N116").Value = -11882.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 545706.7
$ws.Range("I2").Value = 736288.8
$ws.Range("K2").Value = 736288.8
$ws.Range("M2").Value = -736175.8
$ws.Range("H32").Value = 2298.9773
$ws.Range("I32").Value = 2303.675
$ws.Range("K32").Value = 2303.675
$ws.Range("M32").Value = -2016.675
$ws.Range("H45").Value = 4698.8667
$ws.Range("I45").Value = 4851.1924
$ws.Range("K45").Value = 4851.1924
$ws.Range("M45").Value = -4474.1924
$ws.Range("H61").Value = 100002750
$ws.Range("J61").Value = 3555
$ws.Range("L61").Value = 3555
$ws.Range("N61").Value = -3979
$ws.Range("H63").Value = 2146.9167
$ws.Range("I63").Value = 2179.4546
$ws.Range("K63").Value = 2179.4546
$ws.Range("M63").Value = -1493.4546
$ws.Range("H66").Value = 2146.9167
$ws.Range("I66").Value = 2179.4546
$ws.Range("K66").Value = 10897.273
$ws.Range("M66").Value = -7465.273000000001
$ws.Range("H74").Value = 37039910
$ws.Range("I74").Value = 40002636
$ws.Range("J74").Value = 5845.5
$ws.Range("K74").Value = 40002636
$ws.Range("L74").Value = 5845.5
$ws.Range("M74").Value = -40001762
$ws.Range("N74").Value = -7593.5
$ws.Range("H77").Value = 37039910
$ws.Range("I77").Value = 40002636
$ws.Range("J77").Value = 5845.5
$ws.Range("K77").Value = 200013180
$ws.Range("L77").Value = 29227.5
$ws.Range("M77").Value = -200008812
$ws.Range("N77").Value = -37963.5
$ws.Range("H96").Value = 21266.6
$ws.Range("J96").Value = 21266.6
$ws.Range("L96").Value = 21266.6
$ws.Range("N96").Value = -26758.6
$ws.Range("H116").Value = 545706.7
$ws.Range("I116").Value = 736288.8
$ws.Range("K116").Value = 736288.8
$ws.Range("M116").Value = -733994.8
$ws.Range("H132").Value = 2780561.2
$ws.Range("I132").Value = 3451030.2
$ws.Range("J132").Value = 2903.5715
$ws.Range("K132").Value = 10353090.6
$ws.Range("L132").Value = 8710.7145
$ws.Range("M132").Value = -10350560.6
$ws.Range("N132").Value = -13770.7145
$ws.Range("H136").Value = 100002750
$ws.Range("J136").Value = 3555
$ws.Range("L136").Value = 10665
$ws.Range("N136").Value = -15765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 545706.7
$ws.Range("I3").Value = 736288.8
$ws.Range("K3").Value = 736288.8
$ws.Range("M3").Value = -736174.8
$ws.Range("H22").Value = 2070.389
$ws.Range("I22").Value = 2416.7273
$ws.Range("K22").Value = 2416.7273
$ws.Range("M22").Value = -2243.7273
$ws.Range("H107").Value = 50855.715
$ws.Range("I107").Value = 945.8421
$ws.Range("J107").Value = 524999.5
$ws.Range("K107").Value = 945.8421
$ws.Range("L107").Value = 524999.5
$ws.Range("M107").Value = 974.1579
$ws.Range("N107").Value = -528839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7948.364
$ws.Range("I31").Value = 5159.1113
$ws.Range("K31").Value = 5159.1113
$ws.Range("M31").Value = -4864.1113
$ws.Range("H34").Value = 7948.364
$ws.Range("I34").Value = 5159.1113
$ws.Range("K34").Value = 5159.1113
$ws.Range("M34").Value = -4957.1113
$ws.Range("H110").Value = 29999
$ws.Range("J110").Value = 29999
$ws.Range("L110").Value = 29999
$ws.Range("N110").Value = -38179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 168700
$ws.Range("I11").Value = 172741.17
$ws.Range("K11").Value = 518223.51
$ws.Range("M11").Value = -518083.51
$ws.Range("H12").Value = 218.47058
$ws.Range("I12").Value = 135.88889
$ws.Range("K12").Value = 407.66667
$ws.Range("M12").Value = -234.66667
$ws.Range("H97").Value = 659.7143
$ws.Range("I97").Value = 429.5
$ws.Range("J97").Value = 966.6667
$ws.Range("K97").Value = 1288.5
$ws.Range("L97").Value = 2900.0001
$ws.Range("M97").Value = -792.5
$ws.Range("N97").Value = -3892.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1013.5909
$ws.Range("I97").Value = 1090.75
$ws.Range("J97").Value = 807.8333
$ws.Range("K97").Value = 1090.75
$ws.Range("L97").Value = 807.8333
$ws.Range("M97").Value = -594.75
$ws.Range("N97").Value = -1799.8333
$ws.Range("H100").Value = 64999
$ws.Range("J100").Value = 64999
$ws.Range("L100").Value = 64999
$ws.Range("N100").Value = -67163
$ws.Range("H107").Value = 1697.25
$ws.Range("I107").Value = 1911.4286
$ws.Range("J107").Value = 198
$ws.Range("K107").Value = 1911.4286
$ws.Range("L107").Value = 198
$ws.Range("M107").Value = 8.57140000000004
$ws.Range("N107").Value = -4038
$ws.Range("H132").Value = 3208948.8
$ws.Range("J132").Value = 3034.6
$ws.Range("L132").Value = 9103.799999999999
$ws.Range("N132").Value = -14163.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2964.6924
$ws.Range("I46").Value = 2093.7778
$ws.Range("K46").Value = 2093.7778
$ws.Range("M46").Value = -1905.7778
$ws.Range("H68").Value = 2249.2
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 1999
$ws.Range("M68").Value = -1250
$ws.Range("H71").Value = 2249.2
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 9995
$ws.Range("M71").Value = -6251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12163
$ws.Range("J62").Value = 9999
$ws.Range("L62").Value = 9999
$ws.Range("N62").Value = -11247
$ws.Range("H65").Value = 12163
$ws.Range("J65").Value = 9999
$ws.Range("L65").Value = 49995
$ws.Range("N65").Value = -56235
$ws.Range("H81").Value = 3819.4
$ws.Range("I81").Value = 3819.4
$ws.Range("K81").Value = 7638.8
$ws.Range("M81").Value = -6577.8
$ws.Range("H84").Value = 3819.4
$ws.Range("I84").Value = 3819.4
$ws.Range("K84").Value = 38194
$ws.Range("M84").Value = -32890
$ws.Range("H113").Value = 1230.72
$ws.Range("I113").Value = 1227.1666
$ws.Range("K113").Value = 3681.4998
$ws.Range("M113").Value = -1511.4998

Write-Output "Applied all Spriggan Profits market data updates"